$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '28.504.07'
Set-TextCell 'E2' '  -0.12%  '
Set-TextCell 'D3' '1.861.83'
Set-TextCell 'E3' '  +0.80%  '
Set-TextCell 'D4' '1.007'
Set-TextCell 'E4' '  +0.41%  '
Set-TextCell 'D5' '326.34'
Set-TextCell 'E5' '  -2.41%  '
Set-TextCell 'D6' '1.005'
Set-TextCell 'E6' '  +0.34%  '
Set-TextCell 'D7' '0.4645'
Set-TextCell 'E7' '  +0.02%  '
Set-TextCell 'D8' '0.3893'
Set-TextCell 'E8' '  +0.02%  '
Set-TextCell 'D9' '0.07889'
Set-TextCell 'E9' '  +0.07%  '
Set-TextCell 'D10' '0.9728'
Set-TextCell 'E10' '  -0.12%  '
Set-TextCell 'D11' '22.26'
Set-TextCell 'E11' '  +0.48%  '
Set-TextCell 'D12' '1.852.67'
Set-TextCell 'E12' '  +2.10%  '
Set-TextCell 'D13' '5.718'
Set-TextCell 'E13' '  -1.54%  '
Set-TextCell 'D14' '6.935'
Set-TextCell 'E14' '  -0.37%  '
Set-TextCell 'D15' '0.06911'
Set-TextCell 'D16' '88.95'
Set-TextCell 'E16' '  +1.43%  '
Set-TextCell 'D17' '1.007'
Set-TextCell 'E17' '  +0.43%  '
Set-TextCell 'D18' '0.000009999'
Set-TextCell 'E18' '  -0.20%  '
Set-TextCell 'D19' '16.85'
Set-TextCell 'E19' '  -1.04%  '
Set-TextCell 'E20' '  +0.18%  '
Set-TextCell 'D21' '28.494.64'
Set-TextCell 'E21' '  -0.24%  '
Set-TextCell 'D22' '5.329'
Set-TextCell 'E22' '  -0.65%  '
Set-TextCell 'D23' '11.06'
Set-TextCell 'E23' '  -0.75%  '
Set-TextCell 'D24' '2.108'
Set-TextCell 'E24' '  -3.06%  '
Set-TextCell 'D25' '2.087.05'
Set-TextCell 'E25' '  +1.43%  '
Set-TextCell 'D26' '155.05'
Set-TextCell 'E26' '  +1.02%  '
Set-TextCell 'D27' '19.28'
Set-TextCell 'E27' '  -0.53%  '
Set-TextCell 'D28' '5.774'
Set-TextCell 'E28' '  -2.97%  '
Set-TextCell 'D29' '1.989'
Set-TextCell 'E29' '  -0.33%  '
Set-TextCell 'D30' '119.21'
Set-TextCell 'E30' '  +1.42%  '
Set-TextCell 'D31' '0.09309'
Set-TextCell 'E31' '  -0.68%  '
Set-TextCell 'D32' '0.9357'
Set-TextCell 'E32' '  -2.88%  '
Set-TextCell 'D33' '5.315'
Set-TextCell 'E33' '  -0.73%  '
Set-TextCell 'D34' '1.334'
Set-TextCell 'E34' '  -0.53%  '
Set-TextCell 'D35' '3.344'
Set-TextCell 'E35' '  -3.27%  '
Set-TextCell 'D36' '0.05842'
Set-TextCell 'E36' '  -3.70%  '
Set-TextCell 'D37' '0.02119'
Set-TextCell 'E37' '  -3.09%  '
Set-TextCell 'E38' '  -0.99%  '
Set-TextCell 'D39' '7.812'
Set-TextCell 'E39' '  +2.40%  '
Set-TextCell 'D40' '0.5631'
Set-TextCell 'E40' '  -0.85%  '
Set-TextCell 'D41' '9.926'
Set-TextCell 'E41' '  -1.59%  '
Set-TextCell 'D42' '0.1772'
Set-TextCell 'E42' '  -1.21%  '
Set-TextCell 'D43' '0.07355'
Set-TextCell 'E43' '  +4.08%  '
Set-TextCell 'D44' '11.65'
Set-TextCell 'E44' '  -0.89%  '
Set-TextCell 'D45' '0.5304'
Set-TextCell 'E45' '  -0.90%  '
Set-TextCell 'D46' '2.166'
Set-TextCell 'E46' '  -9.07%  '
Set-TextCell 'E47' '  -8.99%  '
Set-TextCell 'D48' '1.845'
Set-TextCell 'E48' '  -1.86%  '
Set-TextCell 'D49' '113.80'
Set-TextCell 'E49' '  +0.51%  '
Set-TextCell 'D50' '2.348'
Set-TextCell 'E50' '  +0.38%  '
Set-TextCell 'D51' '1.006'
Set-TextCell 'E51' '  +0.42%  '
